$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 9 (Ano 2025) with refreshed faturamento figures
$ws.Range("B9").Value = 3061092.22
$ws.Range("C9").Value = 484816.99
$ws.Range("D9").Value = 3545909.21
$ws.Range("E9").Value = 13.67257200586927
$ws.Range("F9").Value = 86.32742799413074
$ws.Range("G9").Value = -53.14485532545871
$ws.Range("H9").Value = -44.72090566973021
$ws.Range("I9").Value = 30720
$ws.Range("J9").Value = 1306
$ws.Range("K9").Value = 32026
$ws.Range("L9").Value = 22087
$ws.Range("M9").Value = 160.5428174944537
$ws.Range("N9").Value = 9.605851396128994
